$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new resale-number row for 2024-01-05 22:34:47
$row = 23

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-05"
$ws.Cells.Item($row, 2).Value = "22:34:47"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"
$ws.Cells.Item($row, 5).Value = 140207
$ws.Cells.Item($row, 6).Value = 142930
$ws.Cells.Item($row, 7).Value = 172448
$ws.Cells.Item($row, 8).Value = 147149
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118296
$ws.Cells.Item($row, 11).Value = 224526
$ws.Cells.Item($row, 12).Value = 248996
$ws.Cells.Item($row, 13).Value = 184801
$ws.Cells.Item($row, 14).Value = 110187
$ws.Cells.Item($row, 15).Value = 40586
$ws.Cells.Item($row, 16).Value = 30823
$ws.Cells.Item($row, 17).Value = 72472
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41933
$ws.Cells.Item($row, 20).Value = -1
